$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8월 1주차")

# --- Row 6 block (2021.08.02 meeting) ---
# Row 6 (권태우): "-" -> "회의록 없음" for D,E,F (C already "회의록 없음")
$ws.Range("D6").Value = "회의록 없음"
$ws.Range("E6").Value = "회의록 없음"
$ws.Range("F6").Value = "회의록 없음"

# Rows 7-9 (김형환, 이혜원, 조현근): "-" -> ditto mark '"'
$ws.Range("C7:F9").Value = '"'

# --- Row 15 block (2021.08.03 meeting) ---
$ws.Range("D15").Value = "회의록 없음"
$ws.Range("E15").Value = "회의록 없음"
$ws.Range("F15").Value = "회의록 없음"

$ws.Range("C16:F18").Value = '"'

# --- Row 24-25 block (2021.08.04 meeting) ---
$ws.Range("C24").Value = "Firebase 연동 전 사용자 및 세탁기 정보를  `nclass 생성 후 데이터 자체저장 작업"
$ws.Range("D24").Value = "1. 예략 버튼 세탁기별 1개에서 `n통합된 버튼 1개로 수정 `n2. Firebase Firestore 활용하여 `nSwift 프로젝트와 연동"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"

$ws.Range("E25").Value = "없음"
$ws.Range("F25").Value = "없음"
